$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet has two header rows (row1 + row2) describing the columns,
# followed by data rows 3-15 and a totals row 16.
# The new layout merges the header into a single row 1, adds new leading
# columns (idx, idx2, Name, Date Start, Date End) and renames the measurement
# headers. Deleting the old row 2 shifts everything up by one row, which is
# exactly what's needed (old row 3 -> new row 2, ..., old row 16 -> new row 15).
$ws.Rows(2).Delete()

# Clear out the stale header fragments (content + formatting) left on row 1
# (old E1/G1/I1/J1/K1).
$ws.Rows(1).Clear()

# Write the new unified header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 pick up the slightly different "applyFont-only" header style used by
# the new layout (same 9pt Arial font as the rest of the sheet).
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# The new layout's selection marks the first data row.
$ws.Range("A2:K2").Select()
